$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.317.02'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.588.67'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '''210.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '1.811.80'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.569.07'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '''64.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '26.321.13'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '''7.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.73%  '
$ws.Range("D20").Value = '''210.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").Value = '''8.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''2.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").Value = '''144.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '''15.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").Value = '1.316.76'
$ws.Range("E34").Value = '  +2.79%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.611'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("E37").Value = '  -0.71%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '''1.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.59%  '
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = '''5.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.41%  '
$ws.Range("D43").Value = '''0.766'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '''2.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("D45").Value = '''62.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = '1.724.46'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").Value = '''87.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  -5.31%  '
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").Value = '''0.0979'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.86%  '
$ws.Range("E51").Value = '  -0.36%  '
